{"js": "// The document was edited so that:\n//  1. In the \"Main Views (SPA Frontend)\" table, the \"Blog post\" row's URI\n//     cell changes from \"/{postTitle}\" to \"/{postId}\" (the frontend route\n//     now keys the blog-post URL by post id, in line with the API which\n//     already used {postId}). Word's internal \"last edit\" bookmark\n//     (_GoBack) ends up sitting right after the inserted text, splitting\n//     the run into \"{postId\" + bookmark + \"}\".\n//  2. In the \"API Resources (Node.js Backend)\" table, the \"Post Comments\"\n//     row's description cell had that same _GoBack bookmark sitting\n//     between \"...Authenticated users\" and the trailing \".\", splitting it\n//     into two runs. Since that is no longer the most-recently-edited\n//     spot, the two runs collapse back into one run of plain text and the\n//     bookmark is gone from there.\n//\n// Both runs involved keep their original run formatting (Arial, italic,\n// lang en-US).\n\nconst RPR_XML =\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Arial\" w:hAnsi=\"Arial\" w:eastAsia=\"Arial\" w:cs=\"Arial\"/>' +\n  '<w:i/>' +\n  '<w:lang w:val=\"en-US\"/>' +\n  '</w:rPr>';\n\nfunction wrapPackage(bodyXml) {\n  return (\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body><w:p>' + bodyXml + '</w:p></w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// --- 1. \"/{postTitle}\" -> \"/{postId\" + _GoBack bookmark + \"}\" -----------\n// The literal \"{postTitle}\" occurs exactly once in the whole document (in\n// the \"Blog post\" row of the Main Views table), so a plain body-wide\n// search is unambiguous.\nconst bodySearch = context.document.body.search(\"{postTitle}\", {\n  matchCase: true,\n});\nbodySearch.load(\"items\");\nawait context.sync();\n\nif (bodySearch.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for {postTitle}, found \" +\n      bodySearch.items.length\n  );\n}\nconst postTitleRange = bodySearch.items[0];\n\nconst newRunsXml =\n  \"<w:r>\" + RPR_XML + \"<w:t>{postId</w:t></w:r>\" +\n  '<w:bookmarkStart w:id=\"1\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"1\"/>' +\n  \"<w:r>\" + RPR_XML + \"<w:t>}</w:t></w:r>\";\n\npostTitleRange.insertOoxml(wrapPackage(newRunsXml), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- 2. Merge \"Comment *Available only for Authenticated users\" + \".\" ---\n//        back into a single run, dropping the _GoBack bookmark that used\n//        to sit between them.\nconst commentSearch = context.document.body.search(\n  \"Comment *Available only for Authenticated users\",\n  { matchCase: true }\n);\ncommentSearch.load(\"items\");\nawait context.sync();\n\nif (commentSearch.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the Comment sentence, found \" +\n      commentSearch.items.length\n  );\n}\n\nconst commentStart = commentSearch.items[0].getRange(\"Start\");\nconst commentParagraph = commentSearch.items[0].paragraphs.getFirst();\nconst paragraphEnd = commentParagraph.getRange(\"End\");\nconst fullCommentRange = commentStart.expandTo(paragraphEnd);\nfullCommentRange.load(\"text\");\nawait context.sync();\n\nconst mergedText = fullCommentRange.text; // \"Comment *Available only for Authenticated users.\"\nconst mergedRunXml =\n  \"<w:r>\" + RPR_XML + \"<w:t>\" + mergedText + \"</w:t></w:r>\";\n\nfullCommentRange.insertOoxml(wrapPackage(mergedRunXml), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document was edited so that:\n#  1. In the \"Main Views (SPA Frontend)\" table, the \"Blog post\" row's URI\n#     cell changes from \"/{postTitle}\" to \"/{postId}\" (the frontend route\n#     now keys the blog-post URL by post id, matching the API which\n#     already used {postId}). Word's internal \"last edit\" bookmark\n#     (_GoBack) ends up sitting right after the inserted text, splitting\n#     the run into \"{postId\" + bookmark + \"}\".\n#  2. In the \"API Resources (Node.js Backend)\" table, the \"Post Comments\"\n#     row's description cell had that same _GoBack bookmark sitting\n#     between \"...Authenticated users\" and the trailing \".\", splitting it\n#     into two runs. Since that is no longer the most-recently-edited\n#     spot, the two runs collapse back into a single run of plain text and\n#     the bookmark disappears from there.\n#\n# _GoBack is a hidden, singleton bookmark that Word itself manages - it\n# cannot be removed with Bookmarks(\"_GoBack\").Delete() (Word silently\n# ignores that), and adding a new one elsewhere just creates a second\n# bookmark of the same name instead of relocating it. So the old one is\n# cleared out by genuinely rewriting the text it sits in (which collapses\n# the runs around it and drops it), and the new one is (re)created with\n# Bookmarks.Add at the desired spot.\n\n$d = $word.ActiveDocument\n\n# --- 1. Drop the old _GoBack bookmark --------------------------------\n# It currently sits between \"...Authenticated users\" and the trailing\n# \".\" in the \"Post Comments\" row description. Re-assigning the range's\n# .Text is a no-op if the replacement string is byte-for-byte identical\n# to the current text (the host skips rewriting), so force a genuine\n# change first and then set the final text back.\n$commentRange = $d.Content\n$commentRange.Find.Execute(\"Comment *Available only for Authenticated users.\") | Out-Null\nif (-not $commentRange.Find.Found) {\n    throw \"Could not find the 'Comment *Available...' sentence.\"\n}\n$commentRange.Text = \"ZZZ_TEMP_PLACEHOLDER_ZZZ\"\n\n$commentRange2 = $d.Content\n$commentRange2.Find.Execute(\"ZZZ_TEMP_PLACEHOLDER_ZZZ\") | Out-Null\nif (-not $commentRange2.Find.Found) {\n    throw \"Could not find the placeholder text while merging the Comment sentence.\"\n}\n$commentRange2.Text = \"Comment *Available only for Authenticated users.\"\n\n# --- 2. \"/{postTitle}\" -> \"/{postId}\" plus a fresh _GoBack bookmark ---\n$titleRange = $d.Content\n$titleRange.Find.Execute(\"{postTitle}\") | Out-Null\nif (-not $titleRange.Find.Found) {\n    throw \"Could not find '{postTitle}'.\"\n}\n$titleRange.Text = \"{postId}\"\n\n$titleRange2 = $d.Content\n$titleRange2.Find.Execute(\"{postId}\") | Out-Null\nif (-not $titleRange2.Find.Found) {\n    throw \"Could not find '{postId}' after replacing '{postTitle}'.\"\n}\n# Collapsed range right after \"{postId\" (7 characters in), i.e. right\n# before the closing \"}\", matching where Word's cursor was left after the\n# edit.\n$bookmarkPoint = $d.Range($titleRange2.Start + 7, $titleRange2.Start + 7)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkPoint)\n"}
